$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '72.527.39'
$ws.Range("E2").Value = '  +0.50%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '4.044.19'
$ws.Range("E3").Value = '  +0.36%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '551.17'
$ws.Range("E5").Value = '  +2.44%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.34'
$ws.Range("E6").Value = '  -0.81%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '4.039.70'
$ws.Range("E7").Value = '  +0.60%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.697'
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.758'
$ws.Range("E10").Value = '  +0.56%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.172'
$ws.Range("E11").Value = '  -0.61%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.17'
$ws.Range("E12").Value = '  +13.18%  '
$ws.Range("E13").Value = '  +0.63%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.97'
$ws.Range("E14").Value = '  +1.65%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.693.62'
$ws.Range("E15").Value = '  +0.98%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.051.62'
$ws.Range("E16").Value = '  +0.89%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.49'
$ws.Range("E17").Value = '  +1.82%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.83'
$ws.Range("E18").Value = '  +1.07%  '
$ws.Range("E19").Value = '  +1.01%  '
$ws.Range("E20").Value = '  -0.69%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.522.46'
$ws.Range("E21").Value = '  +1.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '449.41'
$ws.Range("E22").Value = '  +3.69%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '98.14'
$ws.Range("E23").Value = '  -0.86%  '
$ws.Range("E24").Value = '  -0.50%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.33'
$ws.Range("E25").Value = '  +2.83%  '
$ws.Range("E26").Value = '  +1.50%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.24'
$ws.Range("E27").Value = '  +13.72%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.32'
$ws.Range("E28").Value = '  +1.49%  '
$ws.Range("E29").Value = '  +0.54%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.94'
$ws.Range("E30").Value = '  +1.60%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '37.36'
$ws.Range("E31").Value = '  +0.92%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.91'
$ws.Range("E32").Value = '  +15.62%  '
$ws.Range("E33").Value = '  +3.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '13.68'
$ws.Range("E34").Value = '  +1.51%  '
$ws.Range("B35").Value = 'Bittensor'
$ws.Range("C35").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '686.78'
$ws.Range("E35").Value = '  +0.50%  '
$ws.Range("B36").Value = 'InjectiveProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '49.01'
$ws.Range("E36").Value = '  +15.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '67.37'
$ws.Range("E37").Value = '  +1.68%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.453'
$ws.Range("E38").Value = '  +5.94%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0₃0880'
$ws.Range("E39").Value = '  +5.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.150'
$ws.Range("E40").Value = '  -3.56%  '
$ws.Range("E41").Value = '  -3.44%  '
$ws.Range("E42").Value = '  -1.81%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.23'
$ws.Range("E43").Value = '  +17.35%  '
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("E45").Value = '  +1.77%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("E46").Value = '  +0.28%  '
$ws.Range("E47").Value = '  +0.49%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.71'
$ws.Range("E48").Value = '  +2.41%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.59'
$ws.Range("E49").Value = '  +7.46%  '
$ws.Range("B50").Value = 'ApeXProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.36'
$ws.Range("E50").Value = '  +0.34%  '
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.11'
$ws.Range("E51").Value = '  +2.01%  '
